$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 'knee pads 661'
$ws.Range("A2").Value = 'knee pads ocr'
$ws.Range("A3").Value = 'knee pad inserts for tactical pants'
$ws.Range("A4").Value = 'knee pads jbm'
$ws.Range("A5").Value = 'knee protectors for toddlers'
$ws.Range("A6").Value = 'protec knee pads'
$ws.Range("A7").Value = 'nike thermal compression pants for men'
$ws.Range("A8").Value = 'under armour compression pants youth'
$ws.Range("A9").Value = 'under armour compression tights for men'
$ws.Range("A10").Value = 'capri pants adidas'
$ws.Range("A11").Value = 'capri pants exercise'
$ws.Range("A12").Value = 'capri pants nike'
$ws.Range("A13").Value = 'imucci knee pads'
$ws.Range("A14").Value = 'bb knee pads'
$ws.Range("A15").Value = 'knee pad and helmet'
$ws.Range("A16").Value = 'knee pad basketball kids'
$ws.Range("A17").Value = 'knee pad buttons'
$ws.Range("A18").Value = 'knee pad climbing'
$ws.Range("A19").Value = 'knee pad cover'
$ws.Range("A20").Value = 'knee pad cycling'
$ws.Range("A21").Value = 'knee pad inserts for work pants'
$ws.Range("A22").Value = 'knee pad pants women'
$ws.Range("A23").Value = 'knee pad pink'
$ws.Range("A24").Value = 'knee pad scooter'
$ws.Range("A25").Value = 'knee pad skating'
$ws.Range("A26").Value = 'knee pad wheels'
$ws.Range("A27").Value = 'neoprene knee pads'
$ws.Range("A28").Value = 'skating knee pads'
$ws.Range("A29").Value = 'knee pads capezio'
$ws.Range("A30").Value = 'knee pads canoe'
$ws.Range("A31").Value = 'knee pads cycle'
$ws.Range("A32").Value = 'knee pads firefighter'
$ws.Range("A33").Value = 'knee pads gymnastics'
$ws.Range("A34").Value = 'knee pads multicam'
$ws.Range("A35").Value = 'knee pads protec'
$ws.Range("A36").Value = 'knee pads shin'
$ws.Range("A37").Value = 'knee pads silver'
$ws.Range("A38").Value = 'skate knee pads'
$ws.Range("A39").Value = 'ski knee pad'
$ws.Range("A40").Value = 'kids knee protector'
$ws.Range("A41").Value = 'mens warming compression pants'
$ws.Range("A42").Value = 'mens workout tights'
$ws.Range("A43").Value = 'womens knee pads basketball'
$ws.Range("A44").Value = 'nike pro compression tights men'
$ws.Range("A45").Value = 'mens basketball pants'
$ws.Range("A46").Value = 'basketball knee pads for kids boys'
$ws.Range("A47").Value = 'nike basketball tights'
$ws.Range("A48").Value = 'military pants with knee pads'
$ws.Range("A49").Value = 'asics knee pads'
$ws.Range("A50").Value = 'knee pads xlarge'
$ws.Range("A51").Value = 'motorcycle knee pads men'
$ws.Range("A52").Value = 'gray baseball pants mens'
$ws.Range("A53").Value = 'baseball pants mens knickers'
$ws.Range("A54").Value = 'kids basketball knee pads youth'
$ws.Range("A55").Value = 'youth knee pads basketball for kids'
$ws.Range("A56").Value = 'workout leggings for men'
$ws.Range("A57").Value = 'legging for men nike'
$ws.Range("A58").Value = 'adidas capris men'
$ws.Range("A59").Value = 'youth knee and elbow pads'
$ws.Range("A60").Value = 'knee pads for dancers'
$ws.Range("A61").Value = 'knee pad bathtub'
$ws.Range("A62").Value = 'under armour compression tights men'
$ws.Range("A63").Value = 'knee pad for dancers'
$ws.Range("A64").Value = 'nike youth basketball tights'
$ws.Range("A65").Value = 'compression pants women'
$ws.Range("A66").Value = 'mens workout tights pants'
$ws.Range("A67").Value = 'white nike compression pants men'
$ws.Range("A68").Value = 'defender mens compression pants'
$ws.Range("A69").Value = 'mens nike basketball pants'
$ws.Range("A70").Value = 'nike basketball pants men'
$ws.Range("A71").Value = 'mens leggings compression nike'
$ws.Range("A72").Value = 'white compression pants men'
$ws.Range("A73").Value = 'od green pants with knee pads'
$ws.Range("A74").Value = 'send knee pad'
$ws.Range("A75").Value = 'children knee pads'
$ws.Range("A76").Value = 'smith knee pads'
$ws.Range("A77").Value = 'youth xl football pants'
$ws.Range("A78").Value = 'labor knee pads'
$ws.Range("A79").Value = 'rubber knee pads'
$ws.Range("A80").Value = 'compression pants nike'
$ws.Range("A81").Value = 'athletic capris for women'
$ws.Range("A82").Value = 'skins tights men'
$ws.Range("A83").Value = 'supportive knee pads'
$ws.Range("A84").Value = 'dye knee pads'
$ws.Range("A85").Value = 'dancers knee pads'
$ws.Range("A86").Value = 'apex knee pads'
$ws.Range("A87").Value = 'elbow and knee pads'
$ws.Range("A88").Value = 'fuse knee pads'
$ws.Range("A89").Value = 'ama knee pads'
$ws.Range("A90").Value = 'kp knee pads'
$ws.Range("A91").Value = 'adidas tights men'
$ws.Range("A92").Value = 'caterpillar knee pads'
$ws.Range("A93").Value = 'husky knee pads'
$ws.Range("A94").Value = 'mens leggings white'
$ws.Range("A95").Value = 'nike youth compression pants'
$ws.Range("A96").Value = 'muscle leggings men'
$ws.Range("A97").Value = 'mens compression pants under armour'
$ws.Range("A98").Value = 'airsoft pants with knee pads'
$ws.Range("A99").Value = 'football pants youth with pads'
$ws.Range("A100").Value = 'lotus leggings men'
